$d = $word.ActiveDocument

# Locate the "Latacunga, " run to be replaced with a MERGEFIELD for the
# branch office city, followed by a ", " separator run.
$finder = $d.Content
$found = $finder.Find.Execute("Latacunga, ", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)

if (-not $found) {
    throw "Could not find 'Latacunga, ' text to replace"
}

# Re-materialize a plain Range over the same span: some Find host ranges
# (e.g. $d.Content) behave like an append range under InsertXML, so we
# rebuild a fresh Range(start, end) which really replaces its contents.
$start = $finder.Start
$end = $finder.End
$r = $d.Range($start, $end)

# Clear the run's text first and collapse to an insertion point: when a
# non-empty Range directly abuts a complex field's begin fldChar,
# InsertXML on the un-cleared range can duplicate that neighboring
# fldChar run. Clearing first avoids that.
$r.Text = ""

$rPr = '<w:rPr><w:rFonts w:ascii="Book Antiqua" w:hAnsi="Book Antiqua"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr>'

$newRunsXml = (
    '<w:r w:rsidRPr="00FF3375">' + $rPr + '<w:fldChar w:fldCharType="begin"/></w:r>' +
    '<w:r>' + $rPr + '<w:instrText xml:space="preserve"> MERGEFIELD =consultation.branch_office.city \* MERGEFORMAT </w:instrText></w:r>' +
    '<w:r>' + $rPr + '<w:fldChar w:fldCharType="separate"/></w:r>' +
    '<w:r>' + $rPr + '<w:t>' + [char]0x00AB + '=consultation.branch_office.city' + [char]0x00BB + '</w:t></w:r>' +
    '<w:r>' + $rPr + '<w:fldChar w:fldCharType="end"/></w:r>' +
    '<w:r>' + $rPr + '<w:t xml:space="preserve">, </w:t></w:r>'
)

$xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
    '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
    '<pkg:xmlData>' +
    '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
    '<w:body><w:p>' + $newRunsXml + '</w:p></w:body>' +
    '</w:document>' +
    '</pkg:xmlData></pkg:part></pkg:package>'

$r.InsertXML($xml)

Write-Output "done"
